$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates scraped from the coinranking.com crypto table refresh.
# Values that look like plain numbers are prefixed with a literal leading
# apostrophe (quote-prefix) so Excel stores them as TEXT, not coerced floats
# (e.g. "0.500" must stay "0.500", not become 0.5).

$ws.Range("D2").Value = "57.467.98"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "3.105.74"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'525.42"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'136.78"
$ws.Range("E6").Value = "  -3.60%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.103.74"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").Value = "'0.107"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "'0.394"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "3.637.25"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D15").Value = "'25.24"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").Value = "'0.0000163"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "57.471.44"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "3.098.99"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'5.92"
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").Value = "'12.36"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "'7.85"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").Value = "'345.23"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D24").Value = "'67.60"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'0.500"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "0.0₃0892"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'7.42"
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'1.87"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'6.04"
$ws.Range("E32").Value = "  -6.85%  "
$ws.Range("D33").Value = "'20.79"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").Value = "'4.92"
$ws.Range("E34").Value = "  +6.84%  "
$ws.Range("D35").Value = "'1.15"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "'158.13"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").Value = "'6.05"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'25.85"
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("D39").Value = "'1.23"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "'1.62"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("D41").Value = "'0.0659"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'4.10"
$ws.Range("E42").Value = "  +3.81%  "
$ws.Range("D43").Value = "'0.699"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").Value = "3.142.18"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "2.372.38"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("D46").Value = "'36.62"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("D49").Value = "'0.975"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "'5.97"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").Value = "'19.77"
$ws.Range("E51").Value = "  -3.61%  "
